# Generate Report for Handback
# Refresh the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" values for the 0e249273-... row (row 2) across
# the Overview, zh-cn and de-de sheets, as a handback status report regeneration
# would do.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G2").Value = "2016-08-18 22:50:33"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-08-18 22:50:28"
$zhcn.Range("K2").Value = "2016-08-18 22:50:45"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-08-18 22:50:33"
$dede.Range("K2").Value = "2016-08-18 22:50:52"
